# ToDo.xlsx update: "Finished scala code for AutoMPG Dataset"
#
# Marks every "...Scala" model row as "Done" (they were all "To-Do"
# before), leaving the matching "...Python" rows untouched, and leaves
# the selection where the author last clicked (B21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$doneCells = @("B7", "B9", "B11", "B15", "B17", "B19", "B23", "B25", "B27")
foreach ($addr in $doneCells) {
    $ws.Range($addr).Value = "Done"
}

$null = $ws.Range("B21").Select()
